# Update the hourly HIGH/LOW/CLOSE/LTP/VOL/9:25 CLOSE snapshot values
# on Sheet1 (rows 2-29, columns B-G) with the latest refreshed figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 2538.6
$ws.Range("C2").Value = 2506.1
$ws.Range("D2").Value = 2519
$ws.Range("E2").Value = 2519.3
$ws.Range("F2").Value = 56
$ws.Range("G2").Value = 2516

$ws.Range("B3").Value = 385
$ws.Range("C3").Value = 381.05
$ws.Range("D3").Value = 381.75
$ws.Range("E3").Value = 381.65
$ws.Range("F3").Value = 13
$ws.Range("G3").Value = 384.2

$ws.Range("B4").Value = 1558.8
$ws.Range("C4").Value = 1536.8
$ws.Range("D4").Value = 1540.6
$ws.Range("E4").Value = 1541
$ws.Range("F4").Value = 21
$ws.Range("G4").Value = 1540.55

$ws.Range("B5").Value = 7449.95
$ws.Range("C5").Value = 7370
$ws.Range("D5").Value = 7395
$ws.Range("E5").Value = 7410.65
$ws.Range("F5").Value = 6
$ws.Range("G5").Value = 7419.45

$ws.Range("B6").Value = 245.75
$ws.Range("C6").Value = 238.8
$ws.Range("D6").Value = 243.5
$ws.Range("E6").Value = 243.4
$ws.Range("F6").Value = 152
$ws.Range("G6").Value = 239.4

$ws.Range("B7").Value = 200
$ws.Range("C7").Value = 196.2
$ws.Range("D7").Value = 197.5
$ws.Range("E7").Value = 197.8
$ws.Range("F7").Value = 153
$ws.Range("G7").Value = 197

$ws.Range("B8").Value = 284.45
$ws.Range("C8").Value = 270.9
$ws.Range("D8").Value = 281.25
$ws.Range("E8").Value = 282.1
$ws.Range("F8").Value = 418
$ws.Range("G8").Value = 272.25

$ws.Range("B9").Value = 544.45
$ws.Range("C9").Value = 531.75
$ws.Range("D9").Value = 540.15
$ws.Range("E9").Value = 540
$ws.Range("F9").Value = 80
$ws.Range("G9").Value = 533.25

$ws.Range("B10").Value = 3399.05
$ws.Range("C10").Value = 3366.1
$ws.Range("D10").Value = 3371.7
$ws.Range("E10").Value = 3370.9
$ws.Range("F10").Value = 4
$ws.Range("G10").Value = 3395

$ws.Range("B11").Value = 148.3
$ws.Range("C11").Value = 145.05
$ws.Range("D11").Value = 147.4
$ws.Range("E11").Value = 147.5
$ws.Range("F11").Value = 180
$ws.Range("G11").Value = 145.9

$ws.Range("B12").Value = 1269.55
$ws.Range("C12").Value = 1255.05
$ws.Range("D12").Value = 1262
$ws.Range("E12").Value = 1262.05
$ws.Range("F12").Value = 22
$ws.Range("G12").Value = 1257.65

$ws.Range("B13").Value = 1632.9
$ws.Range("C13").Value = 1608
$ws.Range("D13").Value = 1623.2
$ws.Range("E13").Value = 1623.4
$ws.Range("F13").Value = 172
$ws.Range("G13").Value = 1611.8

$ws.Range("B14").Value = 478.8
$ws.Range("C14").Value = 467.7
$ws.Range("D14").Value = 476.5
$ws.Range("E14").Value = 476.15
$ws.Range("F14").Value = 56
$ws.Range("G14").Value = 471.55

$ws.Range("B15").Value = 975
$ws.Range("C15").Value = 961.25
$ws.Range("D15").Value = 971.8
$ws.Range("E15").Value = 970.55
$ws.Range("F15").Value = 130
$ws.Range("G15").Value = 963.7

$ws.Range("B16").Value = 1459.95
$ws.Range("C16").Value = 1434.6
$ws.Range("D16").Value = 1447.55
$ws.Range("E16").Value = 1448.85
$ws.Range("F16").Value = 28
$ws.Range("G16").Value = 1437.05

$ws.Range("B17").Value = 1474.75
$ws.Range("C17").Value = 1464.7
$ws.Range("D17").Value = 1468.45
$ws.Range("E17").Value = 1469.6
$ws.Range("F17").Value = 31
$ws.Range("G17").Value = 1465.95

$ws.Range("B18").Value = 704.2
$ws.Range("C18").Value = 695.5
$ws.Range("D18").Value = 696.55
$ws.Range("E18").Value = 698.3
$ws.Range("F18").Value = 15
$ws.Range("G18").Value = 703.65

$ws.Range("B19").Value = 457
$ws.Range("C19").Value = 450.2
$ws.Range("D19").Value = 453.4
$ws.Range("E19").Value = 452.95
$ws.Range("F19").Value = 31
$ws.Range("G19").Value = 453

$ws.Range("B20").Value = 1578
$ws.Range("C20").Value = 1565.05
$ws.Range("D20").Value = 1567.3
$ws.Range("E20").Value = 1567.55
$ws.Range("F20").Value = 12
$ws.Range("G20").Value = 1573.25

$ws.Range("B21").Value = 307.45
$ws.Range("C21").Value = 300.8
$ws.Range("D21").Value = 300.9
$ws.Range("E21").Value = 301.65
$ws.Range("F21").Value = 47
$ws.Range("G21").Value = 302.7

$ws.Range("B22").Value = 2456
$ws.Range("C22").Value = 2422.95
$ws.Range("D22").Value = 2451
$ws.Range("E22").Value = 2448.2
$ws.Range("F22").Value = 115
$ws.Range("G22").Value = 2433.5

$ws.Range("B23").Value = 588.5
$ws.Range("C23").Value = 579.7
$ws.Range("D23").Value = 583.4
$ws.Range("E23").Value = 583.45
$ws.Range("F23").Value = 140
$ws.Range("G23").Value = 580.9

$ws.Range("B24").Value = 626.8
$ws.Range("C24").Value = 618.75
$ws.Range("D24").Value = 620.1
$ws.Range("E24").Value = 619.95
$ws.Range("F24").Value = 8
$ws.Range("G24").Value = 622.1

$ws.Range("B25").Value = 1093.7
$ws.Range("C25").Value = 1081.6
$ws.Range("D25").Value = 1085
$ws.Range("E25").Value = 1085.3
$ws.Range("F25").Value = 4
$ws.Range("G25").Value = 1091.4

$ws.Range("B26").Value = 628.9
$ws.Range("C26").Value = 614.85
$ws.Range("D26").Value = 626.95
$ws.Range("E26").Value = 627.25
$ws.Range("F26").Value = 161
$ws.Range("G26").Value = 615.4

$ws.Range("B27").Value = 276.5
$ws.Range("C27").Value = 265.45
$ws.Range("D27").Value = 269.4
$ws.Range("E27").Value = 268.85
$ws.Range("F27").Value = 530
$ws.Range("G27").Value = 266.3

$ws.Range("B28").Value = 131.1
$ws.Range("C28").Value = 129.2
$ws.Range("D28").Value = 129.6
$ws.Range("E28").Value = 129.5
$ws.Range("F28").Value = 343
$ws.Range("G28").Value = 130.75

$ws.Range("B29").Value = 8533.5
$ws.Range("C29").Value = 8400
$ws.Range("D29").Value = 8430
$ws.Range("E29").Value = 8430.549999999999
$ws.Range("F29").Value = 1
$ws.Range("G29").Value = 8532.950000000001
